# Query 1: added % in query directly
# Rename headers from raw flight counts to percentages, and convert the
# row data (on_time/late/cancelled flight counts) into percentages of the
# row total, rounded to 2 decimal places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "on_time_percentage"
$ws.Range("D1").Value = "late_percentage"
$ws.Range("E1").Value = "cancelled_percentage"

# Convert each data row's counts to percentages of the row total
for ($row = 2; $row -le 9; $row++) {
    $onTime = $ws.Cells.Item($row, 3).Value()
    $late = $ws.Cells.Item($row, 4).Value()
    $cancelled = $ws.Cells.Item($row, 5).Value()

    $total = $onTime + $late + $cancelled

    $ws.Cells.Item($row, 3).Value = [Math]::Round(($onTime / $total) * 100, 2)
    $ws.Cells.Item($row, 4).Value = [Math]::Round(($late / $total) * 100, 2)
    $ws.Cells.Item($row, 5).Value = [Math]::Round(($cancelled / $total) * 100, 2)
}
